$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 272.6
$ws.Range("I107").Value = 268.46155
$ws.Range("J107").Value = 299.5
$ws.Range("K107").Value = 268.46155
$ws.Range("L107").Value = 299.5
$ws.Range("M107").Value = 1651.53845
$ws.Range("N107").Value = -4139.5

$ws.Range("H132").Value = 1139.44
$ws.Range("I132").Value = 1095.0526
$ws.Range("J132").Value = 1280
$ws.Range("K132").Value = 3285.1578
$ws.Range("L132").Value = 3840
$ws.Range("M132").Value = -755.1578
$ws.Range("N132").Value = -8900

$ws.Range("I137").Value = 4332.3335
$ws.Range("J137").Value = 4999
$ws.Range("K137").Value = 12997.0005
$ws.Range("L137").Value = 14997
$ws.Range("M137").Value = -10447.0005
$ws.Range("N137").Value = -20097

$ws.Range("H138").Value = 4318.2563
$ws.Range("I138").Value = 1838.909
$ws.Range("J138").Value = 5292.2856
$ws.Range("K138").Value = 5516.727000000001
$ws.Range("L138").Value = 15876.8568
$ws.Range("M138").Value = -376.7270000000008
$ws.Range("N138").Value = -26156.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 124.75
$ws.Range("I5").Value = 100
$ws.Range("J5").Value = 149.5
$ws.Range("K5").Value = 100
$ws.Range("L5").Value = 149.5
$ws.Range("M5").Value = 12
$ws.Range("N5").Value = -373.5

$ws.Range("H32").Value = 3624.8164
$ws.Range("I32").Value = 3387.8333
$ws.Range("J32").Value = 15000
$ws.Range("K32").Value = 3387.8333
$ws.Range("L32").Value = 15000
$ws.Range("M32").Value = -3100.8333

$ws.Range("H61").Value = 1275.75
$ws.Range("I61").Value = 1034.5
$ws.Range("J61").Value = 2240.75
$ws.Range("K61").Value = 1034.5
$ws.Range("L61").Value = 2240.75
$ws.Range("M61").Value = -822.5
$ws.Range("N61").Value = -2664.75

$ws.Range("H74").Value = 4878935.5
$ws.Range("I74").Value = 6451369.5
$ws.Range("J74").Value = 4389.2
$ws.Range("K74").Value = 6451369.5
$ws.Range("L74").Value = 4389.2
$ws.Range("M74").Value = -6450495.5

$ws.Range("H77").Value = 4878935.5
$ws.Range("I77").Value = 6451369.5
$ws.Range("J77").Value = 4389.2
$ws.Range("K77").Value = 32256847.5
$ws.Range("L77").Value = 21946
$ws.Range("M77").Value = -32252479.5

$ws.Range("H97").Value = 1950.1111
$ws.Range("I97").Value = 1129
$ws.Range("J97").Value = 3592.3333
$ws.Range("K97").Value = 1129
$ws.Range("L97").Value = 3592.3333
$ws.Range("M97").Value = -633
$ws.Range("N97").Value = -4584.3333

$ws.Range("H136").Value = 1275.75
$ws.Range("I136").Value = 1034.5
$ws.Range("J136").Value = 2240.75
$ws.Range("K136").Value = 3103.5
$ws.Range("L136").Value = 6722.25
$ws.Range("M136").Value = -553.5
$ws.Range("N136").Value = -11822.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 124.75
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 149.5
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 149.5
$ws.Range("M4").Value = 15
$ws.Range("N4").Value = -379.5

$ws.Range("H134").Value = 3770.5
$ws.Range("I134").Value = 3722.6843
$ws.Range("J134").Value = 4073.3333
$ws.Range("K134").Value = 11168.0529
$ws.Range("L134").Value = 12219.9999
$ws.Range("M134").Value = -8633.052899999999
$ws.Range("N134").Value = -17289.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 228
$ws.Range("I7").Value = 120
$ws.Range("J7").Value = 379.2
$ws.Range("K7").Value = 120
$ws.Range("L7").Value = 379.2
$ws.Range("M7").Value = -7
$ws.Range("N7").Value = -605.2

$ws.Range("H58").Value = 2741.125
$ws.Range("I58").Value = 2577.6
$ws.Range("J58").Value = 3013.6667
$ws.Range("K58").Value = 2577.6
$ws.Range("L58").Value = 3013.6667
$ws.Range("M58").Value = -2374.6
$ws.Range("N58").Value = -3419.6667

$ws.Range("H99").Value = 4437
$ws.Range("I99").Value = 4437
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4437
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -2939
$ws.Range("N99").ClearContents()

$ws.Range("H126").Value = 4437
$ws.Range("I126").Value = 4437
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 13311
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -10841
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 2021.7142
$ws.Range("I132").Value = 1547.6428
$ws.Range("J132").Value = 2969.8572
$ws.Range("K132").Value = 4642.928400000001
$ws.Range("L132").Value = 8909.571599999999
$ws.Range("M132").Value = -2112.928400000001

$ws.Range("H134").Value = 1921.2
$ws.Range("I134").Value = 1963.8422
$ws.Range("J134").Value = 1111
$ws.Range("K134").Value = 5891.5266
$ws.Range("L134").Value = 3333
$ws.Range("M134").Value = -3356.5266
$ws.Range("N134").Value = -8403

$ws.Range("H136").Value = 2741.125
$ws.Range("I136").Value = 2577.6
$ws.Range("J136").Value = 3013.6667
$ws.Range("K136").Value = 7732.799999999999
$ws.Range("L136").Value = 9041.000100000001
$ws.Range("M136").Value = -5182.799999999999
$ws.Range("N136").Value = -14141.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1375.125
$ws.Range("I131").Value = 974.1667
$ws.Range("J131").Value = 1508.7778
$ws.Range("K131").Value = 2922.5001
$ws.Range("L131").Value = 4526.3334
$ws.Range("M131").Value = 2117.4999
$ws.Range("N131").Value = -14606.3334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 15000
$ws.Range("I70").Value = 15000
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 15000
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -14730

$ws.Range("H73").Value = 15000
$ws.Range("I73").Value = 15000
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 15000
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -14064

$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 9000
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -13940

$ws.Range("H132").Value = 1508.1482
$ws.Range("I132").Value = 910.6818
$ws.Range("J132").Value = 4137
$ws.Range("K132").Value = 2732.0454
$ws.Range("L132").Value = 12411
$ws.Range("M132").Value = -202.0454

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1871.5385
$ws.Range("I126").Value = 1683.8
$ws.Range("J126").Value = 2497.3333
$ws.Range("K126").Value = 5051.4
$ws.Range("L126").Value = 7491.999899999999
$ws.Range("M126").Value = -2581.4
$ws.Range("N126").Value = -12431.9999

$ws.Range("H136").Value = 2434.7856
$ws.Range("I136").Value = 1899
$ws.Range("J136").Value = 3399.2
$ws.Range("K136").Value = 5697
$ws.Range("L136").Value = 10197.6
$ws.Range("M136").Value = -3147
$ws.Range("N136").Value = -15297.6
